$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Determine the last used row on the sheet (data starts at row 2).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

# Column C holds the "Förändrad" (changed) date; bump every value in
# rows 2..lastRow from 45179 to 45180 (one day later).
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45179) {
        $cell.Value2 = 45180
    }
}
